$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers for the I0 and IF columns (column I / J), matching the
# bold/bordered header style already used by the other header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Populate I (I0) and J (IF) for each data row (2-37). By default I0 = 1
# and IF mirrors the existing IP (column H) value for that row.
for ($r = 2; $r -le 37; $r++) {
    $ipVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $ipVal
}

# A couple of rows deviate from the default I0=1 / IF=IP pattern.
$ws.Cells.Item(35, 9).Value = 4
$ws.Cells.Item(35, 10).Value = 7

$ws.Cells.Item(37, 9).Value = 3
$ws.Cells.Item(37, 10).Value = 4
